# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Cell B11 on the "Rules" sheet changes from the text "R40" to the text "1".
# Assigning .Value directly would let Excel auto-infer a numeric type for
# "1", so instead enter it as a text formula and convert it to a plain
# value in place (Copy + PasteSpecial values), which preserves the
# existing cell style/format and keeps the result as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false
